$wb = $excel.ActiveWorkbook

$survey = $wb.Worksheets.Item("survey")
$settings = $wb.Worksheets.Item("settings")

# constraint for Survey_Location_Manual: '. = ''' -> '. != '''
$survey.Range("H7").Value = ". != ''"

# fix typo: retailer_cutomer_region -> retailer_customer_region
$survey.Range("B23").Value = "retailer_customer_region"

# fix typo: retailers_market_estimate -> retailer_market_estimate
$survey.Range("B29").Value = "retailer_market_estimate"

# fix typo: reatiler_lpg_rejection -> retailer_lpg_rejection
$survey.Range("B37").Value = "retailer_lpg_rejection"

# bump form __version__ calculation value
# (leading apostrophe doubled to escape Excel's "text prefix" quote semantics)
$survey.Range("J42").Value = "''vSDc6jKQNyUmrTchR9SySA'"

# bump settings version string
$settings.Range("B2").Value = "12 (2022-11-10 09:11:38)"

$wb.Save()
